# "use of DA slide" — append a new "Where is Data Science Needed?" slide
# (Title and Content layout) at the end of the deck.

function Set-BulletArial($rng) {
    $rng.ParagraphFormat.Bullet.Font.Name = "Arial"
    $rng.ParagraphFormat.Bullet.Character = 8226
    $rng.ParagraphFormat.Bullet.Visible = $true
}

$p = $ppt.ActivePresentation

# Title + Content layout (matches the other body slides in this deck).
$s = $p.Slides.Add($p.Slides.Count + 1, 2)

# --- Title -----------------------------------------------------------
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Where is Data Science Needed?"
$title.LanguageID = "en-GB"
$title.Font.Bold = $true

# --- Body bullets ------------------------------------------------------
$body = $s.Shapes.Item(2).TextFrame.TextRange

$body.Text = "To foresee delays for flight/ship/train etc…"
$body.LanguageID = "en-GB"
Set-BulletArial($body)

$r2 = $body.InsertAfter("`rTo predict when a mechanical equipment is about to break")
$r2.LanguageID = "en-GB"
Set-BulletArial($r2)

$r3 = $body.InsertAfter("`rTo create personalized promotional offers")
$r3.LanguageID = "en-GB"
Set-BulletArial($r3)

$r4 = $body.InsertAfter("`rTo find the best suited time to deliver goods / send emails")
$r4.LanguageID = "en-GB"
Set-BulletArial($r4)

$r5 = $body.InsertAfter("`rTo forecast the next years revenue for a company")
$r5.LanguageID = "en-GB"
Set-BulletArial($r5)

$r6 = $body.InsertAfter("`rTo ")
$r6.LanguageID = "en-GB"
Set-BulletArial($r6)
$r6b = $body.InsertAfter("analyze")
$r6b.LanguageID = "en-GB"
$r6c = $body.InsertAfter(" health benefit of treatments")
$r6c.LanguageID = "en-GB"

$r7 = $body.InsertAfter("`rTo predict who will win elections")
$r7.LanguageID = "en-GB"
Set-BulletArial($r7)

$r8 = $body.InsertAfter("`rTo change/update equipment just before they fail")
$r8.LanguageID = "en-GB"
Set-BulletArial($r8)
